$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update names for the two remaining attendee rows
$ws.Range("B2").Value = "DavidLondoño"
$ws.Range("B3").Value = "CarlosRiaño"

# Remove the now-unused rows (rows 4 through 11) entirely so the sheet
# shrinks back down to A1:C3
$ws.Range("A4:C11").EntireRow.Delete()
